$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Moorings")

# Recover Date changed from 42323 (2015-11-15) to 42087 (2015-03-24)
$ws.Range("G2").Value = 42087

# Notes column: add "Glider lost"
$ws.Range("L2").Value = "Glider lost"

# Update selection to G7 on the Moorings sheet
$ws.Activate()
$ws.Range("G7").Select()
